$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextCell 2 4 "28.303.46"
Set-TextCell 2 5 "  -2.41%  "

Set-TextCell 3 4 "1.869.17"
Set-TextCell 3 5 "  -2.09%  "

Set-TextCell 4 5 "  +0.06%  "

Set-TextCell 5 4 "318.61"
Set-TextCell 5 5 "  -2.00%  "

Set-TextCell 6 4 "1.003"
Set-TextCell 6 5 "  +0.10%  "

Set-TextCell 7 4 "0.4405"
Set-TextCell 7 5 "  -4.12%  "

Set-TextCell 8 4 "0.3698"
Set-TextCell 8 5 "  -3.37%  "

Set-TextCell 9 4 "0.07500"
Set-TextCell 9 5 "  -2.83%  "

Set-TextCell 10 4 "0.9375"
Set-TextCell 10 5 "  -4.37%  "

Set-TextCell 11 4 "21.43"
Set-TextCell 11 5 "  -2.93%  "

Set-TextCell 12 4 "1.899.94"
Set-TextCell 12 5 "  +0.94%  "

Set-TextCell 13 4 "6.704"
Set-TextCell 13 5 "  -3.34%  "

Set-TextCell 14 4 "5.465"
Set-TextCell 14 5 "  -3.64%  "

Set-TextCell 15 4 "0.06885"
Set-TextCell 15 5 "  -2.15%  "

Set-TextCell 16 4 "1.004"
Set-TextCell 16 5 "  -0.01%  "

Set-TextCell 17 4 "82.02"
Set-TextCell 17 5 "  -2.10%  "

Set-TextCell 18 4 "0.000009035"
Set-TextCell 18 5 "  -4.55%  "

Set-TextCell 19 4 "1.003"
Set-TextCell 19 5 "  +0.06%  "

Set-TextCell 20 4 "15.90"
Set-TextCell 20 5 "  -4.63%  "

Set-TextCell 21 4 "28.301.56"
Set-TextCell 21 5 "  -2.27%  "

Set-TextCell 22 4 "5.117"
Set-TextCell 22 5 "  -3.80%  "

Set-TextCell 23 4 "10.84"
Set-TextCell 23 5 "  -0.42%  "

Set-TextCell 24 4 "2.093.98"
Set-TextCell 24 5 "  -1.35%  "

Set-TextCell 25 4 "2.028"
Set-TextCell 25 5 "  -3.12%  "

Set-TextCell 26 5 "  -2.14%  "

Set-TextCell 27 4 "18.41"
Set-TextCell 27 5 "  -3.62%  "

Set-TextCell 28 4 "5.318"
Set-TextCell 28 5 "  -6.02%  "

Set-TextCell 29 4 "113.48"
Set-TextCell 29 5 "  -3.43%  "

Set-TextCell 30 4 "1.720"
Set-TextCell 30 5 "  -7.25%  "

Set-TextCell 31 4 "0.09032"
Set-TextCell 31 5 "  -2.79%  "

Set-TextCell 32 4 "0.7977"
Set-TextCell 32 5 "  -7.53%  "

Set-TextCell 33 4 "4.859"
Set-TextCell 33 5 "  -4.23%  "

Set-TextCell 34 4 "1.174"
Set-TextCell 34 5 "  -5.70%  "

Set-TextCell 35 4 "2.935"
Set-TextCell 35 5 "  -2.68%  "

Set-TextCell 36 2 "Frax"
Set-TextCell 36 3 "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextCell 36 4 "1.003"
Set-TextCell 36 5 "  +0.09%  "

Set-TextCell 37 2 "TrustWalletToken"
Set-TextCell 37 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell 37 4 "1.130"
Set-TextCell 37 5 "  -2.05%  "

Set-TextCell 38 4 "0.05447"
Set-TextCell 38 5 "  -4.91%  "

Set-TextCell 39 2 "VeChain"
Set-TextCell 39 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell 39 4 "0.01969"
Set-TextCell 39 5 "  -3.30%  "

Set-TextCell 40 2 "MXToken"
Set-TextCell 40 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell 40 4 "3.018"
Set-TextCell 40 5 "  +6.13%  "

Set-TextCell 41 4 "7.108"
Set-TextCell 41 5 "  -4.02%  "

Set-TextCell 42 4 "0.5261"
Set-TextCell 42 5 "  -4.54%  "

Set-TextCell 43 5 "  -4.01%  "

Set-TextCell 44 4 "8.731"
Set-TextCell 44 5 "  -6.54%  "

Set-TextCell 45 4 "0.06758"
Set-TextCell 45 5 "  -1.02%  "

Set-TextCell 46 4 "0.4889"
Set-TextCell 46 5 "  -5.75%  "

Set-TextCell 47 5 "  -4.89%  "

Set-TextCell 48 4 "107.66"
Set-TextCell 48 5 "  -2.92%  "

Set-TextCell 49 4 "1.947"
Set-TextCell 49 5 "  -4.99%  "

Set-TextCell 50 2 "PaxDollar"
Set-TextCell 50 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell 50 4 "1.002"
Set-TextCell 50 5 "  +0.10%  "

Set-TextCell 51 2 "NEARProtocol"
Set-TextCell 51 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell 51 4 "1.675"
Set-TextCell 51 5 "  -5.92%  "
